$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.9
$ws.Range("I2").Value = 4.2
$ws.Range("J2").Value = 2.6
$ws.Range("L2").Value = 4.75
$ws.Range("Q2").Value = 2.08
$ws.Range("R2").Value = 1.73
$ws.Range("Y2").Value = 8.5
$ws.Range("Z2").Value = 15
$ws.Range("AC2").Value = 9
$ws.Range("AI2").Value = 21
$ws.Range("AO2").Value = 10
$ws.Range("AP2").Value = 21
$ws.Range("AQ2").Value = 34
$ws.Range("AW2").Value = 6
$ws.Range("G3").Value = 2.38
$ws.Range("I3").Value = 3.6
$ws.Range("J3").Value = 3.4
$ws.Range("L3").Value = 4.5
$ws.Range("AM3").Value = 51
$ws.Range("G4").Value = 1.62
$ws.Range("I4").Value = 6.25
$ws.Range("J4").Value = 2.25
$ws.Range("M4").Value = 1.08
$ws.Range("N4").Value = 8
$ws.Range("X4").Value = 7
$ws.Range("Y4").Value = 9
$ws.Range("AC4").Value = 8
$ws.Range("AE4").Value = 19
$ws.Range("AG4").Value = 1250
$ws.Range("AN4").Value = 3.5
$ws.Range("AQ4").Value = 29
$ws.Range("AU4").Value = 9
$ws.Range("G6").Value = 1.8
$ws.Range("I6").Value = 4.75
$ws.Range("J6").Value = 2.5
$ws.Range("L6").Value = 5
$ws.Range("Q6").Value = 2.15
$ws.Range("R6").Value = 1.67
$ws.Range("AB6").Value = 29
$ws.Range("AE6").Value = 17
$ws.Range("AG6").Value = 800
$ws.Range("AI6").Value = 21
$ws.Range("AJ6").Value = 15
$ws.Range("AM6").Value = 41
$ws.Range("AO6").Value = 9.5
$ws.Range("AS6").Value = 151
$ws.Range("AU6").Value = 8.5
$ws.Range("AW6").Value = 6
$ws.Range("AX6").Value = 26
$ws.Range("AY6").Value = 34
$ws.Range("BA6").Value = 126
$ws.Range("I8").Value = 2.87
$ws.Range("J8").Value = 2.9
$ws.Range("L8").Value = 3.5
$ws.Range("V8").Value = 2.25
$ws.Range("W8").Value = 9.75
$ws.Range("Y8").Value = 9
$ws.Range("AA8").Value = 17
$ws.Range("AG8").Value = 300
$ws.Range("AK8").Value = 35
$ws.Range("AP8").Value = 19
$ws.Range("AS8").Value = 250
$ws.Range("AW8").Value = 4.9
$ws.Range("AX8").Value = 16
$ws.Range("AZ8").Value = 75
$ws.Range("G9").Value = 2.1
$ws.Range("H9").Value = 3.5
$ws.Range("J9").Value = 2.65
$ws.Range("K9").Value = 2.18
$ws.Range("L9").Value = 3.7
$ws.Range("S9").Value = 1.36
$ws.Range("V9").Value = 2.18
$ws.Range("Z9").Value = 20
$ws.Range("AX9").Value = 17.5
$ws.Range("G10").Value = 2.77
$ws.Range("H10").Value = 3.15
$ws.Range("I10").Value = 2.5
$ws.Range("J10").Value = 3.35
$ws.Range("O10").Value = 1.29
$ws.Range("W10").Value = 9
$ws.Range("X10").Value = 14.5
$ws.Range("AB10").Value = 29
$ws.Range("AD10").Value = 6.1
$ws.Range("AH10").Value = 8.75
$ws.Range("AK10").Value = 28
$ws.Range("AM10").Value = 27
$ws.Range("AO10").Value = 15
$ws.Range("AP10").Value = 22
$ws.Range("AR10").Value = 100
$ws.Range("AS10").Value = 250
$ws.Range("AU10").Value = 6.7
$ws.Range("AV10").Value = 55
$ws.Range("AX10").Value = 13.5
$ws.Range("I11").Value = 6.5
$ws.Range("J11").Value = 1.98
$ws.Range("K11").Value = 2.25
$ws.Range("M11").Value = 1.05
$ws.Range("N11").Value = 7.6
$ws.Range("O11").Value = 1.28
$ws.Range("Q11").Value = 1.85
$ws.Range("R11").Value = 1.88
$ws.Range("T11").Value = 2.75
$ws.Range("W11").Value = 6.2
$ws.Range("AB11").Value = 30
$ws.Range("AC11").Value = 7.6
$ws.Range("AM11").Value = 80
$ws.Range("AN11").Value = 3.2
$ws.Range("AT11").Value = 2.75
$ws.Range("AY11").Value = 40
